$d = $word.ActiveDocument

$old = "nombres y apellido del paciente, nº de documento, fecha de nacimiento, estado civil, nacionalidad, lugar de nacimiento, provincia, obra social, religión, nivel de estudios, trabajo o profesión, domicilio actual, teléfono, celular, médico de cabecera, teléfono del médico de cabecera, si posee servicios de emergencia, cuál de ellos, si el paciente se encuentra privado de la libertad y donde."
$new = "nombres y apellido del paciente, tipo y nº de documento, fecha de nacimiento, estado civil, nacionalidad, lugar de nacimiento, provincia, obra social, religión, nivel de estudios, trabajo o profesión, calle, numero de calle, barrio, ciudad, domicilio anterior, teléfono, celular, médico de cabecera, teléfono del médico de cabecera, si posee servicios de emergencia, cuál de ellos, si el paciente se encuentra privado de la libertad y donde."

$count = 0
$found = $true
while ($found) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) { $count = $count + 1 }
}
Write-Output "Replaced: $count"
